$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("Sheet4")
$ws1 = $wb.Worksheets.Item("Sheet1")

$longStr1 = 'dffffffffffffffffff fh uhurhguhgh u  htri hg  hri ri hg  I hgi hri hrgir hgihrighsrigrs '
$sdd = 'sdd'
$dd = 'dd'
$longStr2 = 'fdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie h fdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie hfdf jfij ji  ij  jij  jijfijf jfijfi j jdfithe u ie h'
$sOnly = 's'

# Sheet4 tab (file xl/worksheets/sheet3.xml): add rows 19-24 with new shared strings,
# in the same order they first appear in the diff so shared-string indices line up.
$ws4.Range("A19").Value = $longStr1
$ws4.Range("A20").Value = $sdd
$ws4.Range("A21").Value = $sdd
$ws4.Range("A22").Value = $dd
$ws4.Range("A23").Value = $sdd
$ws4.Range("A24").Value = $longStr1

# Sheet1 tab (file xl/worksheets/sheet1.xml): add row 4 and row 10
$ws1.Range("A4").Value = $longStr2
$ws1.Range("A10").Value = $sOnly

# Update selections + active sheet/tab.
# Sheet4 tab loses tabSelected, gets a single-cell selection at A24.
$ws4.Activate() | Out-Null
$ws4.Range("A24").Select() | Out-Null

# Sheet1 tab becomes the active tab, with selection at F12.
$ws1.Activate() | Out-Null
$ws1.Range("F12").Select() | Out-Null
